$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for C, D, H (B and A, E, F, G stay unchanged).
# Note: this runtime's ColumnWidth setter round-trips through a pixel
# conversion on save that adds a constant +5/6 offset, so subtract it
# here to land exactly on the intended integer character-widths.
$widthFudge = 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 49 - $widthFudge
$ws.Columns.Item(4).ColumnWidth = 49 - $widthFudge
$ws.Columns.Item(8).ColumnWidth = 38 - $widthFudge

# Opportunity IDs in column A are text (e.g. "1330536"), not numbers -
# force text storage so numeric-looking ids don't get reinterpreted as numbers.
$ws.Range("A2:A12").NumberFormat = "@"

# New data for rows 2-12 (A:H)
$data = @(
    @("1331015", "https://aiesec.org/opportunity/global-talent/1331015", "Marketing", "Sincan, İstasyon, 06934 Sincan/Ankara, Türkiye", "No", "14 applicants", "9 - 12 Weeks", "Maksan Lift"),
    @("1331005", "https://aiesec.org/opportunity/global-talent/1331005", "Industrial Engineer", "Istanbul, İstanbul, Türkiye", "No", "15 applicants", "6 - 18 Months", "Unixpadel"),
    @("1329657", "https://aiesec.org/opportunity/global-talent/1329657", "Sales and Marketing Specialist", "Ankara, Türkiye", "No", "49 applicants", "9 - 12 Weeks", "Metaform"),
    @("1328272", "https://aiesec.org/opportunity/global-talent/1328272", "Graphic Designer", "Tunis, Tunisia", "No", "20 applicants", "9 - 12 Weeks", "La fabrique"),
    @("1328206", "https://aiesec.org/opportunity/global-talent/1328206", "Power BI Specialist", "Frankfurt am Main, Deutschland", "No", "187 applicants", "3 - 6 Months", "Greyfood GmbH"),
    @("1327779", "https://aiesec.org/opportunity/global-talent/1327779", "Business Development & Foreign Trade Assistant", "Sincan, İstasyon, 06934 Sincan/Ankara, Türkiye", "No", "58 applicants", "9 - 12 Weeks", "ATM Atılım Teknik Makine"),
    @("1327156", "https://aiesec.org/opportunity/global-talent/1327156", "Security Master Domain Intern", "Heerlen, Netherlands", "No", "93 applicants", "6 - 18 Months", "APG Heerlen"),
    @("1321256", "https://aiesec.org/opportunity/global-talent/1321256", "Businesses Development Intern", "Lahore, Punjab, Pakistan", "No", "21 applicants", "9 - 12 Weeks", "Devsinc."),
    @("1320899", "https://aiesec.org/opportunity/global-talent/1320899", "F&B Operations Manager", "Hong Kong", "No", "78 applicants", "6 - 18 Months", "Treehouse"),
    @("1308313", "https://aiesec.org/opportunity/global-talent/1308313", "Chemical Engineering", "Ankara, Türkiye", "No", "222 applicants", "6 - 18 Months", "AKA SU TEKNOLOJİLERİ ANONİM ŞİRKETİ"),
    @("1301473", "https://aiesec.org/opportunity/global-talent/1301473", "Digital Marketing Specialist in Istanbul", "İstanbul, Türkiye", "No", "33 applicants", "6 - 18 Months", "Raff Textile")
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
